# datosVerificacionAlgoritmo.xlsx - "Se modifica secuencial y se agrega funcion (libreria)"
#
# Adds a small 5x5 "image" verification block (Q5:U9) and a 3x3 "mask"
# verification block (M7:O9) next to the existing algorithm-trace table,
# plus a tiny "SUM" library-function result cell (L12/M12), and moves the
# active selection to K18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new cellXfs style: like the existing left/center-aligned style (s=1)
# but with the thin border already used by the header row (s=2), minus the
# bold font. Build it by copying the header cell's format (font+border+
# alignment) and then turning the bold font back off, instead of composing
# borders by hand (which would mint a brand-new border record).
$ws.Range("G5").Copy() | Out-Null
$ws.Range("Q5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Q5").Font.Bold = $false

# Propagate that same format (border + left/center alignment, regular font)
# across every cell of the two new data blocks before filling values in.
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("M7:O9").PasteSpecial(-4122) | Out-Null
$ws.Range("Q5:U9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 5x5 "image" block: Q5:U9
$ws.Range("Q5:U9").Value = @(
    @(6, 5, 6, 3, 4),
    @(4, 5, 6, 3, 2),
    @(2, 2, 5, 1, 1),
    @(2, 5, 9, 3, 2),
    @(1, 5, 6, 3, 2)
)

# --- 3x3 "mask" block: M7:O9
$ws.Range("M7:O9").Value = @(
    @(2, 3, 4),
    @(5, 1, 5),
    @(6, 2, 6)
)

# --- "SUM" library-function result, next to the algorithm trace table
$ws.Range("L12").Value = "SUM"
$ws.Range("M12").Value = 0
$ws.Range("L12:M12").Font.Bold = $false

# Match the column widths the workbook ships with for the new helper
# columns flanking the two blocks (narrow spacer columns around the data).
$narrowWidth = 2 - 0.8333333333333334
foreach ($col in 13..15) { $ws.Columns.Item($col).ColumnWidth = $narrowWidth }
foreach ($col in 17..21) { $ws.Columns.Item($col).ColumnWidth = $narrowWidth }

# Move the active selection as recorded after the edit.
$ws.Range("K18").Select() | Out-Null
